$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.451
$ws.Range("C2").Value = 0.601

$ws.Range("B3").Value = 0.301
$ws.Range("C3").Value = 0.451

$ws.Range("B4").Value = 0.151
$ws.Range("C4").Value = 0.301
